$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the single "Date" column into three separate columns: Year, Month, Day.
# Insert two new blank columns before column A, shifting the old A:S header
# row (and its per-cell styles) right to C:U.
$ws.Range("A:B").EntireColumn.Insert() | Out-Null

# Populate the two freshly-inserted columns.
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Month"

# Old column A ("Date") is now column C - relabel it "Day" (its original
# date-style formatting is kept as-is).
$ws.Range("C1").Value = "Day"

# Give the new Year/Month header cells the same formatting as the other
# plain text headers (e.g. the shifted "River" header now in D1), instead of
# default/no style.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null

# Match the workbook's recorded selection after the edit.
$ws.Range("C2").Select() | Out-Null
